# EventReminderMain: Emailing Capability Implemented
#
# The only meaningful content change between before/after is on Sheet2:
# cell A2 held "country " (with a trailing space) and is corrected to
# "country" (no trailing space) so it matches the string already used on
# Sheet1!A1. That also lets Excel collapse the now-unused "country "
# shared-string entry when it re-saves the workbook.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Fix the trailing-space typo in the country column header value.
$ws2.Range("A2").Value = "country"

# Reflect the cursor/selection landing on A2 after the edit, as captured
# in the saved sheetView for Sheet2.
$ws2.Activate()
$ws2.Range("A2").Select()
